# Insert a new column before column D (pushes existing D:K data to E:L)
# and populate the new column D with the latest fiscal year figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; everything from D to K (and beyond) shifts right by one.
$ws.Columns("D").Insert()

# Copy the (now shifted) column E formatting into the new column D, but only
# for the rows that actually carry data in that block (skip the section
# title rows 5, 6, 37 and 79, which have no D:K cells at all).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Match the new column's width to its neighbours (same bestFit look as D:H).
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# New fiscal-year data for the newly inserted column D.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 167939600
$ws.Range("D9").Value = 163291400
$ws.Range("D10").Value = 4648300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 302900
$ws.Range("D15").Value = 465100
$ws.Range("D17").Value = 166519700
$ws.Range("D18").Value = 1419900
$ws.Range("D20").Value = -52900
$ws.Range("D21").Value = 1877200
$ws.Range("D22").Value = 189600
$ws.Range("D23").Value = 1177400
$ws.Range("D24").Value = 174100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1003300
$ws.Range("D27").Value = 1045800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 612600
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 52900
$ws.Range("D33").Value = 1658400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1658400

$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 2492500
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 11314200
$ws.Range("D44").Value = 11918500
$ws.Range("D45").Value = 169100
$ws.Range("D46").Value = 25894400
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1892400
$ws.Range("D49").Value = 9612100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 270900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 37669800
$ws.Range("D57").Value = 26836900
$ws.Range("D58").Value = 151700
$ws.Range("D59").Value = 881200
$ws.Range("D60").Value = 27869700
$ws.Range("D61").Value = 4510800
$ws.Range("D62").Value = 2239400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 34737000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 3720600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2932800
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 1658400
$ws.Range("D83").Value = 510100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1411400
$ws.Range("D91").Value = -336400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1111100
$ws.Range("D96").Value = -333000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -242900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 57400
